{"js": "// Append the new \"State (3)\" section at the end of the document body,\n// right after the last existing paragraph (\"ConcreteSubscriber ...\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the very last paragraph in the document.\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfunction addParagraph(text, opts) {\n  opts = opts || {};\n  const p = anchor.insertParagraph(text || \"\", Word.InsertLocation.after);\n  // Always set a style explicitly (rather than leaving it alone) \u2014 Word\n  // clones the *entire* pPr (style + numbering + alignment) from the\n  // previous paragraph otherwise, and none of the newly authored\n  // paragraphs below are meant to inherit ConcreteSubscriber's bullet /\n  // \"justify both\" formatting.\n  if (opts.heading2) {\n    p.styleBuiltIn = Word.Style.heading2;\n  } else if (opts.listParagraph) {\n    p.styleBuiltIn = Word.Style.listParagraph;\n    p.attachToList(2, 0); // reuse the existing numId=2 bulleted list\n  } else {\n    p.styleBuiltIn = Word.Style.normal;\n  }\n  // None of the new paragraphs in the target carry the body's usual\n  // \"justify both\" direct formatting, so make sure it stays left/default.\n  p.alignment = Word.Alignment.left;\n  anchor = p;\n  return p;\n}\n\nfunction addRun(paragraph, text) {\n  return paragraph.insertText(text, Word.InsertLocation.end);\n}\n\n// 1) Heading2 \u2014 \"State (3)\"\naddParagraph(\"State (3)\", { heading2: true });\n\n// 2) Intro paragraph (built from several runs of identical formatting \u2014\n//    they collapse into the same text when concatenated).\nlet p2 = addParagraph(\n  \"State patern zamenjuje state mashinu (finite state machine) koja je implementirana preko \"\n);\naddRun(p2, \"mnogo \");\naddRun(p2, \"uslovnih opreatora. \");\naddRun(p2, \"Obicno su stanja state masine samo gomila fildova sa drugim vrednostima\");\naddRun(\n  p2,\n  \". Najveca mana ovakvih state mashina je da broj i kompleksnost uslova drasticno raste kako se dodaju novi uslovi a pogotovo novi bihejvori u zavisnosti od stanja.\"\n);\n\n// 3) Context paragraph (+ the sentence that, in the original edit, fell\n//    right after a page break).\nlet p3 = addParagraph(\n  \"State patern nam to resava tako sto kreira nove klase za sva moguca stanja i ponasanja nekog objekta. Umesto da implementira sva ponasanja i vrednosti za sva stanja, imamo objekt Context koji samo cuva referencu na konkretan state objekat koji predstavlja njegovo trenutno stanje i delegira sav posao na taj objekat. Da bi se promenilo stanje konteksta u novo samo se zameni referenca na aktivan state objekat. \"\n);\naddRun(\n  p3,\n  \"Ovo je naravno moguce samo ako sve klase koje predstavljaju stanja implementiraju jedan isti interfejs i sam kontekst sa objektima state klasa komunicira sa njima kroz taj interfejs.\"\n);\n\n// 4) Strategy comparison paragraph.\naddParagraph(\n  \"Ova struktura je slicna Strategy paternu, ali glavna razlika je sto u State paternu stanja mogu da znaju jedna za druge i iniciraju tranziciju iz jednog u drugo stanje, dok strategije ne znaju jedna za drugu.\"\n);\n\n// 5) \"Klase koje ucestvuju:\" lead-in.\naddParagraph(\"Klase koje ucestvuju:\");\n\n// 6-8) Bulleted list items (ListParagraph, numId 2).\nlet p6 = addParagraph(\"Context \", { listParagraph: true });\naddRun(\n  p6,\n  \"\u2013 definise interfejs preko koga se klijenti obracaju stanjima, cuva instancu na ConcreteState podklasu koja definise trenutno stanje\"\n);\n\naddParagraph(\"State \u2013 definise interfejs koji enkapsulira ponasanje klase za odredjeno stanje (State)\", {\n  listParagraph: true,\n});\n\naddParagraph(\n  \"ConcreteState \u2013 svaka klasa implementira interfejs State i konkretno ponasanje vezano za stanje konteksta.\",\n  { listParagraph: true }\n);\n\nawait context.sync();\n", "ps1": "# Append the new \"State (3)\" section at the end of the document, right\n# after the last existing paragraph (\"ConcreteSubscriber ...\").\n\n$d = $word.ActiveDocument\n\nfunction Add-BodyParagraph([string]$text) {\n    $endRange = $d.Content\n    $endRange.Collapse(0)            # wdCollapseEnd\n    $endRange.InsertParagraphAfter()\n    $endRange.Collapse(0)\n    $p = $d.Paragraphs.Last\n    # Drop any inherited bullet numbering before touching style/text so no\n    # stray <w:numPr numId=\"0\"/> is left behind.\n    $p.Range.ListFormat.RemoveNumbers()\n    $p.Style = \"Normal\"\n    $p.Alignment = 0                 # wdAlignParagraphLeft (no w:jc)\n    $p.Range.Text = $text\n    return $p\n}\n\nfunction Add-Heading2Paragraph([string]$text) {\n    $endRange = $d.Content\n    $endRange.Collapse(0)\n    $endRange.InsertParagraphAfter()\n    $endRange.Collapse(0)\n    $p = $d.Paragraphs.Last\n    $p.Range.ListFormat.RemoveNumbers()\n    $p.Style = \"Heading 2\"\n    $p.Alignment = 0\n    $p.Range.Text = $text\n    return $p\n}\n\nfunction Add-ListParagraph([string]$text) {\n    $endRange = $d.Content\n    $endRange.Collapse(0)\n    $endRange.InsertParagraphAfter()\n    $endRange.Collapse(0)\n    $p = $d.Paragraphs.Last\n    $p.Style = \"List Paragraph\"\n    $p.Alignment = 0\n    $p.Range.Text = $text\n    # Reuse the existing numId=2 bulleted list already used elsewhere in\n    # the document (ilvl 0).\n    $p.Range.ListFormat.ApplyListTemplateWithLevel($d.ListTemplates.Item(2), $false, 0, $false, $false)\n    return $p\n}\n\nfunction Add-Run($paragraph, [string]$text) {\n    $r = $paragraph.Range\n    $r.Collapse(0)\n    $r.InsertAfter($text)\n}\n\n# 1) Heading2 -- \"State (3)\"\nAdd-Heading2Paragraph \"State (3)\" | Out-Null\n\n# 2) Intro paragraph, built from several runs of identical formatting.\n$p2 = Add-BodyParagraph \"State patern zamenjuje state mashinu (finite state machine) koja je implementirana preko \"\nAdd-Run $p2 \"mnogo \"\nAdd-Run $p2 \"uslovnih opreatora. \"\nAdd-Run $p2 \"Obicno su stanja state masine samo gomila fildova sa drugim vrednostima\"\nAdd-Run $p2 \". Najveca mana ovakvih state mashina je da broj i kompleksnost uslova drasticno raste kako se dodaju novi uslovi a pogotovo novi bihejvori u zavisnosti od stanja.\"\n\n# 3) Context paragraph (+ sentence that originally fell after a page break).\n$p3 = Add-BodyParagraph \"State patern nam to resava tako sto kreira nove klase za sva moguca stanja i ponasanja nekog objekta. Umesto da implementira sva ponasanja i vrednosti za sva stanja, imamo objekt Context koji samo cuva referencu na konkretan state objekat koji predstavlja njegovo trenutno stanje i delegira sav posao na taj objekat. Da bi se promenilo stanje konteksta u novo samo se zameni referenca na aktivan state objekat. \"\nAdd-Run $p3 \"Ovo je naravno moguce samo ako sve klase koje predstavljaju stanja implementiraju jedan isti interfejs i sam kontekst sa objektima state klasa komunicira sa njima kroz taj interfejs.\"\n\n# 4) Strategy comparison paragraph.\nAdd-BodyParagraph \"Ova struktura je slicna Strategy paternu, ali glavna razlika je sto u State paternu stanja mogu da znaju jedna za druge i iniciraju tranziciju iz jednog u drugo stanje, dok strategije ne znaju jedna za drugu.\" | Out-Null\n\n# 5) \"Klase koje ucestvuju:\" lead-in.\nAdd-BodyParagraph \"Klase koje ucestvuju:\" | Out-Null\n\n# 6-8) Bulleted list items (List Paragraph, numId 2).\n$p6 = Add-ListParagraph \"Context \"\nAdd-Run $p6 \"\u2013 definise interfejs preko koga se klijenti obracaju stanjima, cuva instancu na ConcreteState podklasu koja definise trenutno stanje\"\n\nAdd-ListParagraph \"State \u2013 definise interfejs koji enkapsulira ponasanje klase za odredjeno stanje (State)\" | Out-Null\n\nAdd-ListParagraph \"ConcreteState \u2013 svaka klasa implementira interfejs State i konkretno ponasanje vezano za stanje konteksta.\" | Out-Null\n\nWrite-Output \"done\"\n"}
